$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Comment for the commits field (E9), entered before the profile info so the
# shared-string table keeps the original author's insertion order.
$ws.Range("E9").Value = "commits are not much but are valuable and tested"

# SoftUni student info (GitHub username / full name)
$ws.Range("C4").Value = "yavor2000"
$ws.Range("C5").Value = "Явор Митев"

# Admin home page added -> answer the Admin Home Screen question
$ws.Range("C34").Value = "Yes"

# Updated commit counts
$ws.Range("C8").Value = 11
$ws.Range("C9").Value = 20

# Newly filled basic-option scores
$ws.Range("C26").Value = 10
$ws.Range("C27").Value = 5
$ws.Range("C29").Value = 5
$ws.Range("C30").Value = 5

# Restore the saved view/selection state from the edit session
$ws.Application.ActiveWindow.ScrollRow = 13
[void]$ws.Range("C34").Select()
